$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the total "VALOR MORA" amount
$ws.Range("E11").Value = 160000

# Update the "Cant. Periodos" count (only one period remains now)
$ws.Range("F13").Value = 1

# Update the remaining period row: period label + accrued value
$ws.Range("E16").Value = "2508"
$ws.Range("G16").Value = 4000000

# Remove the second (now obsolete) period row entirely, shifting rows below up
$ws.Rows("17").Delete()
